# Replace the "Glen Waverley" / "Melbourne" exposure-site rows with an
# updated, much larger list of public exposure sites.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell content (but keep formatting, e.g. the bold /
# bordered header row) so the shared-string table is rebuilt cleanly from
# the new data (rather than retaining now-unused strings like
# "Glen Waverley").
$ws.Cells.ClearContents()

$data = @(
    @("Location",          "Site",                                                                                  "Exposure period",              "Notes",                                                                           "Exist"),
    @("Broadmeadows",       "Craigieburn Line train",                                                                "1.25pm - 1.59pm  9/02/2012",   "Case caught train from Broadmeadows Railway Station to Glenroy Railway Station", "new"),
    @("Coburg",             "Function venue  426 Sydney Rd  Coburg VIC 3058",                                        "7:14pm  11:30pm  6/02/2021",   "Case attended venue",                                                             "new"),
    @("Glenroy",            "513 Eltham to Glenroy bus route  Glenroy Railway Station towards Eltham",               "1.35pm  2.17pm  9/02/2021",    "Case caught bus from Glenroy Railway Station towards Eltham",                    "new"),
    @("Hoppers Crossing",   "Caltex Woolworths  50 Old Geelong Rd  Hoppers Crossing, VIC 3029",                      "6.40am - 7.15am  8/02/21",     "Case attended venue",                                                             "new"),
    @("Hoppers Crossing",   "Coates Hire Werribee  148A Geelong Rd  Hoppers Crossing, VIC 3029",                     "6.45am - 7.30am  8/02/21",     "Case attended venue",                                                             "new"),
    @("Melbourne",          "901 Frankston to Melbourne Airport bus route  Melbourne Airport to Broadmeadows Railway Station", "1:02pm  1:49pm  9/2/2021", "Case caught but from Melbourne Airport to Broadmeadows Railway Station", "new"),
    @("Melbourne",          "Brunetti: Terminal 4, Melbourne Airport",                                               "4:45am - 1:15pm  9/2/2021",    "Case attended venue",                                                             "new"),
    @("Melbourne",          "Brunetti: Terminal 4, Melbourne Airport",                                               "4:45am - 1:15pm, 9/2/2021",    "Case attended venue",                                                             "old"),
    @("South Melbourne",    "Stowe Australia  67  69 Buckhurst St  South Melbourne VIC 3205",                        "10.30am - 10.45am 8/02/2021",  $null,                                                                             "new")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r + 1, $c + 1).Value = $val
        }
    }
}

# Widen the columns to fit the new, longer text.
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(2).ColumnWidth = 59
$ws.Columns.Item(3).ColumnWidth = 24.5
$ws.Columns.Item(4).ColumnWidth = 65.5

# Print setup.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("C7").Select()
